# Apply changes to the insurance (保險) and debt (債務) sheets:
#  - sheet "保險" (insurance, Worksheets.Item(6)): relabel the header row with
#    proper field names, fix the B/C column values that had been shifted by
#    one row, and append the standard trailing metadata columns
#    (property_category/category/date/legislator_name/legislator_id/
#    source_file/index).
#  - sheet "債務" (debt, Worksheets.Item(7)): same idea - relabel the header
#    row, fix the shifted B/C/D/F/G values, and append the same trailing
#    metadata columns plus the "debt" category marker.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 6: 保險 (insurance)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Extend formatting (border/alignment styles) from the existing B:D columns
# out to the new E:K columns, for both the header row and the data rows.
$ws6.Range("B1:D10").Copy()
$ws6.Range("E1:K10").PasteSpecial(-4122)

# Header row
$ws6.Range("B1").Value = "company"
$ws6.Range("C1").Value = "name"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "property_category"
$ws6.Range("F1").Value = "category"
$ws6.Range("G1").Value = "date"
$ws6.Range("H1").Value = "legislator_name"
$ws6.Range("I1").Value = "legislator_id"
$ws6.Range("J1").Value = "source_file"
$ws6.Range("K1").Value = "index"

$insuranceRows = @(
    @{ Row = 2;  Index = 136; Company = "新光人壽";     Name = "千禧金黉貝";            Owner = "郭綺雯" },
    @{ Row = 3;  Index = 137; Company = "新光人壽";     Name = "千禧傳家寶";            Owner = "郭綺雯" },
    @{ Row = 4;  Index = 138; Company = "新光人壽";     Name = "千禧傳家寶";            Owner = "郭綺雯" },
    @{ Row = 5;  Index = 139; Company = "新光人壽";     Name = "寶順養老";              Owner = "郭綺雯" },
    @{ Row = 6;  Index = 140; Company = "新光人壽";     Name = "新住院醫療";            Owner = "郭綺雯" },
    @{ Row = 7;  Index = 141; Company = "新光人壽";     Name = "寶順養老";              Owner = "郭綺雯" },
    @{ Row = 8;  Index = 142; Company = "國際紐約人壽"; Name = "儲蓄壽險";              Owner = "郭綺雯" },
    @{ Row = 9;  Index = 143; Company = "富邦人壽";     Name = "GPLB二十年繳費终生壽險"; Owner = "紀國棟" },
    @{ Row = 10; Index = 144; Company = "中國人壽";     Name = "保誠美滿養老保險";      Owner = "紀國棟" }
)

foreach ($r in $insuranceRows) {
    $row = $r.Row
    $ws6.Range("B$row").Value = $r.Company
    $ws6.Range("C$row").Value = $r.Name
    $ws6.Range("D$row").Value = $r.Owner
    $ws6.Range("E$row").Value = "insurance"
    $ws6.Range("F$row").Value = "normal"
    $ws6.Range("G$row").Value = "2012-04-16"
    $ws6.Range("H$row").Value = "紀國棟"
    $ws6.Range("I$row").Value = 918
    $ws6.Range("J$row").Value = "tmpf6b91"
    $ws6.Range("K$row").Value = $r.Index
}

$ws6.Range("A1:K1").EntireRow.RowHeight = $ws6.Range("A1:K1").EntireRow.RowHeight

# ---------------------------------------------------------------------
# Sheet 7: 債務 (debt)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Extend formatting out to the new H:N columns as well.
$ws7.Range("B1:G5").Copy()
$ws7.Range("H1:N5").PasteSpecial(-4122)

# Header row
$ws7.Range("B1").Value = "species"
$ws7.Range("C1").Value = "debtor"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "total"
$ws7.Range("F1").Value = "register_date"
$ws7.Range("G1").Value = "register_reason"
$ws7.Range("H1").Value = "property_category"
$ws7.Range("I1").Value = "category"
$ws7.Range("J1").Value = "date"
$ws7.Range("K1").Value = "legislator_name"
$ws7.Range("L1").Value = "legislator_id"
$ws7.Range("M1").Value = "source_file"
$ws7.Range("N1").Value = "index"

$debtRows = @(
    @{ Row = 2; Index = 154; Species = "抵押貸款"; Debtor = "紀國棟"; Owner = "臺中市大肚鄉農會臺中市大肚區沙田路";               Total = 1426133;  RegisterDate = "96年03月15日";  RegisterReason = "付房屋貸款" },
    @{ Row = 3; Index = 155; Species = "抵押貸款"; Debtor = "郭綺雯"; Owner = "中國信託商業銀行文心分行臺中市北屯區文心路";     Total = 2914607;  RegisterDate = "96年09月21日";  RegisterReason = "購置房屋" },
    @{ Row = 4; Index = 156; Species = "抵押貸款"; Debtor = "郭綺雯"; Owner = "國泰世華商業銀行豐原分行臺中市豐原區三民路";     Total = 11540590; RegisterDate = "100年05月06日"; RegisterReason = "購置房屋" },
    @{ Row = 5; Index = 157; Species = "抵押貸款"; Debtor = "紀國棟"; Owner = "台新國際商業銀行文心分行臺中市北屯區文心路";     Total = 15673437; RegisterDate = "100年10月12日"; RegisterReason = "購置房屋" }
)

foreach ($r in $debtRows) {
    $row = $r.Row
    $ws7.Range("B$row").Value = $r.Species
    $ws7.Range("C$row").Value = $r.Debtor
    $ws7.Range("D$row").Value = $r.Owner
    $ws7.Range("E$row").Value = $r.Total
    $ws7.Range("F$row").Value = $r.RegisterDate
    $ws7.Range("G$row").Value = $r.RegisterReason
    $ws7.Range("H$row").Value = "debt"
    $ws7.Range("I$row").Value = "normal"
    $ws7.Range("J$row").Value = "2012-04-16"
    $ws7.Range("K$row").Value = "紀國棟"
    $ws7.Range("L$row").Value = 918
    $ws7.Range("M$row").Value = "tmpf6b91"
    $ws7.Range("N$row").Value = $r.Index
}
